# edit.ps1 - apply Horse_presentation.pptx content updates
# Mirrors the authored diff: titles de-prefixed ("Slide N: " removed / reworded),
# body paragraphs reworded, and source lines replaced with article URLs.
$p = $ppt.ActivePresentation

function Set-RunText($shape, [string]$newText) {
    $tr = $shape.TextFrame.TextRange
    # Reset to an unrelated placeholder first so the host's run-diffing
    # logic can't splice the new text onto a stale run via a shared prefix;
    # this keeps the paragraph down to a single run, matching the source deck.
    $tr.Text = 'placeholder reset value zzz'
    $tr.Text = $newText
}

function Set-SourceText($shape, [string]$newText) {
    # The source textbox carries a leading empty paragraph before the
    # actual citation line; a leading "`r" in .Text keeps that paragraph
    # intact instead of collapsing the shape down to one paragraph.
    $origHeight = $shape.Height
    $tr = $shape.TextFrame.TextRange
    $tr.Text = "`rplaceholder reset value zzz"
    $tr.Text = "`r$newText"
    $shape.Height = $origHeight
}

# Slide 10
$s = $p.Slides.Item(10)
Set-RunText $s.Shapes.Item(1) 'Horses in Modern Society'
Set-RunText $s.Shapes.Item(2) 'In modern society, horses continue to play important roles in areas such as agriculture, therapy, sports, and recreation. They are valued for their beauty, grace, and companionship.'
Set-SourceText $s.Shapes.Item(3) '- https://horse-canada.com/horses-and-history/horses-use-today/'

# Slide 11
$s = $p.Slides.Item(11)
Set-RunText $s.Shapes.Item(1) 'Conclusion'
Set-RunText $s.Shapes.Item(2) 'Horses are fascinating creatures with a rich history and a significant presence in human culture. By understanding and respecting these magnificent animals, we can build strong and meaningful connections with them.'
Set-SourceText $s.Shapes.Item(3) '- https://www.equinesciencenews.com/the-beauty-and-intricacy-of-the-horse/'

# Slide 2
$s = $p.Slides.Item(2)
Set-RunText $s.Shapes.Item(1) 'Introduction to Horses'
Set-RunText $s.Shapes.Item(2) 'Horses are majestic and powerful animals that have been domesticated for thousands of years. They have played a crucial role in human history, from transportation to agriculture to sports.'
Set-SourceText $s.Shapes.Item(3) '- https://en.wikipedia.org/wiki/Horse'

# Slide 3
$s = $p.Slides.Item(3)
Set-RunText $s.Shapes.Item(1) 'Types of Horses'
Set-RunText $s.Shapes.Item(2) 'There are over 300 different breeds of horses, each with its own unique characteristics and abilities. Some common types include Arabian, Thoroughbred, and Quarter Horse.'
Set-SourceText $s.Shapes.Item(3) '- https://www.thesprucepets.com/horse-breeds-1118558'

# Slide 4
$s = $p.Slides.Item(4)
Set-RunText $s.Shapes.Item(1) 'Anatomy of a Horse'
Set-RunText $s.Shapes.Item(2) 'Horses have a complex and intricate anatomy, with strong muscles, a powerful heart, and a digestion system that is unique among animals. Understanding their anatomy is crucial for proper care and handling.'
Set-SourceText $s.Shapes.Item(3) '- https://en.wikipedia.org/wiki/Equine_anatomy'

# Slide 5
$s = $p.Slides.Item(5)
Set-RunText $s.Shapes.Item(1) 'Life Cycle of a Horse'
Set-RunText $s.Shapes.Item(2) 'Horses typically live for 25 to 30 years, with their life cycle consisting of various stages such as foal, yearling, and adult. Proper nutrition and healthcare are essential for a horse''s well-being.'
Set-SourceText $s.Shapes.Item(3) '- https://www.msdvetmanual.com/management-and-nutrition/husbandry-of-horses'

# Slide 6
$s = $p.Slides.Item(6)
Set-RunText $s.Shapes.Item(1) 'Horse Behavior'
Set-RunText $s.Shapes.Item(2) 'Horses are social animals that exhibit a wide range of behaviors, including communication through body language, grooming each other, and forming strong bonds with their herd mates.'
Set-SourceText $s.Shapes.Item(3) '- https://www.horsemagazine.com/thm/2011/08/equine-behaviour-constant-themes-and-new-insights/'

# Slide 7
$s = $p.Slides.Item(7)
Set-RunText $s.Shapes.Item(1) 'Horse Care and Maintenance'
Set-RunText $s.Shapes.Item(2) 'Proper care and maintenance of horses involve providing them with a balanced diet, regular exercise, grooming, and veterinary care. Good horsemanship is essential for keeping horses healthy and happy.'
Set-SourceText $s.Shapes.Item(3) '- https://www.equisearch.com/discoverhorses/horse-care-10-care-tips-horse-care-25281'

# Slide 8
$s = $p.Slides.Item(8)
Set-RunText $s.Shapes.Item(1) 'Horse Riding and Training'
Set-RunText $s.Shapes.Item(2) 'Horse riding and training require skill, patience, and mutual trust between the horse and the rider. Different disciplines such as dressage, show jumping, and western riding offer a variety of ways to enjoy working with horses.'
Set-SourceText $s.Shapes.Item(3) '- https://practicalhorsemanmag.com/training/english-disciplines-explained-14986'

# Slide 9
$s = $p.Slides.Item(9)
Set-RunText $s.Shapes.Item(1) 'Famous Horses in History'
Set-RunText $s.Shapes.Item(2) 'Throughout history, there have been many famous horses that have left a lasting impact on human society. Examples include Bucephalus, Seabiscuit, and Secretariat.'
Set-SourceText $s.Shapes.Item(3) '- https://www.horseandman.com/people-and-places/10-historical-horses-with-amazing-stories/'
